$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list: updated prices / 1h % changes, plus two
# leaderboard swaps (ARBITRUM <-> ImmutableX, NEARProtocol <-> PaxDollar).

# Row 2
$ws.Range("D2").Value = "'27.844.80"
$ws.Range("E2").Value = '  +1.55%  '

# Row 3
$ws.Range("D3").Value = "'1.867.91"
$ws.Range("E3").Value = '  +1.65%  '

# Row 4
$ws.Range("E4").Value = '  +0.66%  '

# Row 5
$ws.Range("D5").Value = "'324.25"
$ws.Range("E5").Value = '  +2.00%  '

# Row 6
$ws.Range("D6").Value = "'1.029"
$ws.Range("E6").Value = '  +0.51%  '

# Row 7
$ws.Range("D7").Value = "'0.4416"
$ws.Range("E7").Value = '  +1.25%  '

# Row 8
$ws.Range("D8").Value = "'0.3819"
$ws.Range("E8").Value = '  +2.57%  '

# Row 9
$ws.Range("D9").Value = "'0.07461"
$ws.Range("E9").Value = '  +1.34%  '

# Row 10
$ws.Range("D10").Value = "'0.8890"
$ws.Range("E10").Value = '  +1.84%  '

# Row 11
$ws.Range("D11").Value = "'21.71"
$ws.Range("E11").Value = '  +1.90%  '

# Row 12
$ws.Range("D12").Value = "'1.880.82"
$ws.Range("E12").Value = '  +1.79%  '

# Row 13
$ws.Range("D13").Value = "'5.574"
$ws.Range("E13").Value = '  +2.03%  '

# Row 14
$ws.Range("D14").Value = "'6.775"
$ws.Range("E14").Value = '  +1.29%  '

# Row 15
$ws.Range("D15").Value = "'0.07193"
$ws.Range("E15").Value = '  +1.18%  '

# Row 16
$ws.Range("D16").Value = "'85.75"
$ws.Range("E16").Value = '  +4.15%  '

# Row 17
$ws.Range("E17").Value = '  +0.70%  '

# Row 18
$ws.Range("D18").Value = "'0.000009142"
$ws.Range("E18").Value = '  +1.64%  '

# Row 19
$ws.Range("E19").Value = '  +0.49%  '

# Row 20
$ws.Range("D20").Value = "'15.59"
$ws.Range("E20").Value = '  +1.30%  '

# Row 21
$ws.Range("D21").Value = "'27.849.59"
$ws.Range("E21").Value = '  +1.52%  '

# Row 22
$ws.Range("D22").Value = "'5.323"
$ws.Range("E22").Value = '  +1.74%  '

# Row 23
$ws.Range("D23").Value = "'11.31"
$ws.Range("E23").Value = '  +1.48%  '

# Row 24
$ws.Range("D24").Value = "'2.112.29"
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
$ws.Range("D25").Value = "'2.036"
$ws.Range("E25").Value = '  +7.64%  '

# Row 26
$ws.Range("D26").Value = "'158.53"
$ws.Range("E26").Value = '  +1.08%  '

# Row 27
$ws.Range("D27").Value = "'18.84"
$ws.Range("E27").Value = '  +1.41%  '

# Row 28
$ws.Range("D28").Value = "'5.409"
$ws.Range("E28").Value = '  +3.42%  '

# Row 29
$ws.Range("D29").Value = "'1.996"
$ws.Range("E29").Value = '  +4.10%  '

# Row 30
$ws.Range("D30").Value = "'118.21"
$ws.Range("E30").Value = '  +2.09%  '

# Row 31
$ws.Range("D31").Value = "'0.09054"
$ws.Range("E31").Value = '  +0.17%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'0.7860"
$ws.Range("E32").Value = '  +3.71%  '

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = "'1.227"
$ws.Range("E33").Value = '  +2.21%  '

# Row 34
$ws.Range("D34").Value = "'3.032"
$ws.Range("E34").Value = '  +5.84%  '

# Row 35
$ws.Range("D35").Value = "'4.602"
$ws.Range("E35").Value = '  +3.00%  '

# Row 36
$ws.Range("D36").Value = "'1.032"
$ws.Range("E36").Value = '  +0.59%  '

# Row 37
$ws.Range("D37").Value = "'1.149"
$ws.Range("E37").Value = '  +0.13%  '

# Row 38
$ws.Range("D38").Value = "'0.01988"
$ws.Range("E38").Value = '  +1.37%  '

# Row 39
$ws.Range("D39").Value = "'0.05332"
$ws.Range("E39").Value = '  +1.66%  '

# Row 40
$ws.Range("D40").Value = "'2.877"
$ws.Range("E40").Value = '  +3.28%  '

# Row 41
$ws.Range("D41").Value = "'0.5229"
$ws.Range("E41").Value = '  +1.46%  '

# Row 42
$ws.Range("D42").Value = "'0.1691"
$ws.Range("E42").Value = '  +1.86%  '

# Row 43
$ws.Range("D43").Value = "'6.922"
$ws.Range("E43").Value = '  +5.54%  '

# Row 44
$ws.Range("D44").Value = "'8.889"
$ws.Range("E44").Value = '  +4.65%  '

# Row 45
$ws.Range("D45").Value = "'111.23"
$ws.Range("E45").Value = '  +2.42%  '

# Row 46
$ws.Range("D46").Value = "'10.81"
$ws.Range("E46").Value = '  +2.41%  '

# Row 47
$ws.Range("D47").Value = "'0.06619"
$ws.Range("E47").Value = '  +4.98%  '

# Row 48
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = "'1.033"
$ws.Range("E48").Value = '  +0.60%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'1.723"
$ws.Range("E49").Value = '  +2.70%  '

# Row 50
$ws.Range("D50").Value = "'0.4748"
$ws.Range("E50").Value = '  +2.75%  '

# Row 51
$ws.Range("D51").Value = "'1.925"
$ws.Range("E51").Value = '  +1.81%  '
